$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 102; this shifts existing rows 102-111 down to 103-112,
# keeping all of their data/styles intact.
$ws.Rows.Item(102).Insert()

# Populate the new row 102 with the values for the added record.
$ws.Cells.Item(102, 1).Value = 10
$ws.Cells.Item(102, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(102, 3).Value = "La Araucanía"
$ws.Cells.Item(102, 4).Value = Get-Date -Year 2022 -Month 8 -Day 10 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(102, 5).Value = 9
$ws.Cells.Item(102, 6).Value = 100114002
$ws.Cells.Item(102, 7).Value = "Camote"
$ws.Cells.Item(102, 8).Value = "Sin especificar"
$ws.Cells.Item(102, 9).Value = "Primera"
$ws.Cells.Item(102, 10).Value = 30
$ws.Cells.Item(102, 11).Value = 20000
$ws.Cells.Item(102, 12).Value = 20000
$ws.Cells.Item(102, 13).Value = 20000
$ws.Cells.Item(102, 14).Value = "$/malla 20 kilos"
$ws.Cells.Item(102, 15).Value = "Perú"
$ws.Cells.Item(102, 16).Value = 1000
$ws.Cells.Item(102, 17).Value = 20
$ws.Cells.Item(102, 18).Value = "Hortaliza"
